# Applies scheduled-runner price/profit updates to the Leve tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 101
$ws.Range("H101").Value = 2547.9
$ws.Range("I101").Value = 1275.7858
$ws.Range("K101").Value = 3827.3574
$ws.Range("M101").Value = -2205.3574
# ALC row 134
$ws.Range("H134").Value = 33328
$ws.Range("J134").Value = 33328
$ws.Range("L134").Value = 33328
$ws.Range("N134").Value = -43468
# ALC row 138
$ws.Range("H138").Value = 2850.791
$ws.Range("I138").Value = 3013
$ws.Range("J138").Value = 2760.2559
$ws.Range("K138").Value = 9039
$ws.Range("L138").Value = 8280.7677
$ws.Range("M138").Value = -3899
$ws.Range("N138").Value = -18560.7677

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 9355.166999999999
$ws.Range("I32").Value = 13847.25
$ws.Range("K32").Value = 13847.25
$ws.Range("M32").Value = -13560.25
# ARM row 45
$ws.Range("H45").Value = 1412.8334
$ws.Range("I45").Value = 695.7143
$ws.Range("J45").Value = 1869.1818
$ws.Range("K45").Value = 695.7143
$ws.Range("L45").Value = 1869.1818
$ws.Range("M45").Value = -318.7143
$ws.Range("N45").Value = -2623.1818
# ARM row 61
$ws.Range("H61").Value = 4171.5674
$ws.Range("I61").Value = 1375.0294
$ws.Range("K61").Value = 1375.0294
$ws.Range("M61").Value = -1163.0294
# ARM row 110
$ws.Range("H110").Value = 1436.4286
$ws.Range("I110").Value = 1342.75
$ws.Range("K110").Value = 1342.75
$ws.Range("M110").Value = 702.25
# ARM row 136
$ws.Range("H136").Value = 4171.5674
$ws.Range("I136").Value = 1375.0294
$ws.Range("K136").Value = 4125.0882
$ws.Range("M136").Value = -1575.0882

$ws = $wb.Worksheets.Item("BSM")
# BSM row 99
$ws.Range("H99").Value = 2921.55
$ws.Range("I99").Value = 2613.5881
$ws.Range("J99").Value = 4666.6665
$ws.Range("K99").Value = 2613.5881
$ws.Range("L99").Value = 4666.6665
$ws.Range("M99").Value = -1115.5881
$ws.Range("N99").Value = -7662.6665
# BSM row 107
$ws.Range("H107").Value = 5099.75
$ws.Range("I107").Value = 5099.75
$ws.Range("K107").Value = 5099.75
$ws.Range("M107").Value = -3179.75

$ws = $wb.Worksheets.Item("CRP")
# CRP row 4
$ws.Range("H4").Value = 29999
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 29999
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 29999
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = -30223
# CRP row 62
$ws.Range("H62").Value = 7376
$ws.Range("J62").Value = 7678.4
$ws.Range("L62").Value = 7678.4
$ws.Range("N62").Value = -8926.4
# CRP row 65
$ws.Range("H65").Value = 7376
$ws.Range("J65").Value = 7678.4
$ws.Range("L65").Value = 38392
$ws.Range("N65").Value = -44632
# CRP row 107
$ws.Range("H107").Value = 263.25
$ws.Range("J107").Value = 237.2
$ws.Range("L107").Value = 237.2
$ws.Range("N107").Value = -4077.2
# CRP row 132
$ws.Range("H132").Value = 2958.7874
$ws.Range("I132").Value = 2927.9768
$ws.Range("K132").Value = 8783.930399999999
$ws.Range("M132").Value = -6253.930399999999
# CRP row 134
$ws.Range("H134").Value = 3754.1025
$ws.Range("I134").Value = 3587.8484
$ws.Range("J134").Value = 4668.5
$ws.Range("K134").Value = 10763.5452
$ws.Range("L134").Value = 14005.5
$ws.Range("M134").Value = -8228.5452
$ws.Range("N134").Value = -19075.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4
$ws.Range("H4").Value = 72994430
$ws.Range("I4").Value = 127613650
$ws.Range("J4").Value = 168816.67
$ws.Range("K4").Value = 382840950
$ws.Range("L4").Value = 506450.01
$ws.Range("M4").Value = -382840838
$ws.Range("N4").Value = -506674.01
# CUL row 5
$ws.Range("H5").Value = 227.25
$ws.Range("I5").Value = 206.5
$ws.Range("J5").Value = 248
$ws.Range("K5").Value = 619.5
$ws.Range("L5").Value = 744
$ws.Range("M5").Value = -507.5
$ws.Range("N5").Value = -968
# CUL row 14
$ws.Range("H14").Value = 775
$ws.Range("I14").Value = 775
$ws.Range("K14").Value = 2325
$ws.Range("M14").Value = -2152
# CUL row 18
$ws.Range("H18").Value = 556.6667
$ws.Range("I18").Value = 111
$ws.Range("J18").Value = 645.8
$ws.Range("K18").Value = 333
$ws.Range("L18").Value = 1937.4
$ws.Range("M18").Value = -164
$ws.Range("N18").Value = -2275.4
# CUL row 23
$ws.Range("H23").Value = 532.25
$ws.Range("I23").Value = 349
$ws.Range("J23").Value = 615.5454999999999
$ws.Range("K23").Value = 1047
$ws.Range("L23").Value = 1846.6365
$ws.Range("M23").Value = -812
$ws.Range("N23").Value = -2316.6365
# CUL row 33
$ws.Range("H33").Value = 5225
$ws.Range("I33").Value = 10000
$ws.Range("K33").Value = 60000
$ws.Range("M33").Value = -59717
# CUL row 109
$ws.Range("H109").Value = 2058.2222
$ws.Range("I109").Value = 2203.125
$ws.Range("K109").Value = 6609.375
$ws.Range("M109").Value = -5569.375
# CUL row 121
$ws.Range("H121").Value = 22311676
$ws.Range("J121").Value = 30415586
$ws.Range("L121").Value = 91246758
$ws.Range("N121").Value = -91249378
# CUL row 135
$ws.Range("H135").Value = 227.25
$ws.Range("I135").Value = 206.5
$ws.Range("J135").Value = 248
$ws.Range("K135").Value = 1858.5
$ws.Range("L135").Value = 2232
$ws.Range("M135").Value = 676.5
$ws.Range("N135").Value = -7302

$ws = $wb.Worksheets.Item("GSM")
# GSM row 5
$ws.Range("H5").Value = 9999
$ws.Range("J5").Value = 9999
$ws.Range("L5").Value = 9999
$ws.Range("N5").Value = -10223
# GSM row 97
$ws.Range("H97").Value = 907.88
$ws.Range("I97").Value = 295.3
$ws.Range("K97").Value = 295.3
$ws.Range("M97").Value = 200.7
# GSM row 126
$ws.Range("H126").Value = 4432.55
$ws.Range("I126").Value = 3951.3635
$ws.Range("J126").Value = 5020.6665
$ws.Range("K126").Value = 11854.0905
$ws.Range("L126").Value = 15061.9995
$ws.Range("M126").Value = -9384.0905
$ws.Range("N126").Value = -20001.9995
# GSM row 132
$ws.Range("H132").Value = 2243.6177
$ws.Range("I132").Value = 2172.6191
$ws.Range("J132").Value = 2358.3076
$ws.Range("K132").Value = 6517.8573
$ws.Range("L132").Value = 7074.9228
$ws.Range("M132").Value = -3987.8573
$ws.Range("N132").Value = -12134.9228
# GSM row 136
$ws.Range("H136").Value = 29800.75
$ws.Range("J136").Value = 29800.75
$ws.Range("L136").Value = 89402.25
$ws.Range("N136").Value = -94502.25

$ws = $wb.Worksheets.Item("LTW")
# LTW row 2
$ws.Range("H2").Value = 24197.2
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 29996.5
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 29996.5
$ws.Range("M2").Value = -888
$ws.Range("N2").Value = -30220.5
# LTW row 22
$ws.Range("H22").Value = 1575.75
$ws.Range("I22").Value = 878.25
$ws.Range("J22").Value = 1924.5
$ws.Range("K22").Value = 878.25
$ws.Range("L22").Value = 1924.5
$ws.Range("M22").Value = -583.25
$ws.Range("N22").Value = -2514.5
# LTW row 27
$ws.Range("H27").Value = 1575.75
$ws.Range("I27").Value = 878.25
$ws.Range("J27").Value = 1924.5
$ws.Range("K27").Value = 878.25
$ws.Range("L27").Value = 1924.5
$ws.Range("M27").Value = -771.25
$ws.Range("N27").Value = -2138.5
# LTW row 46
$ws.Range("H46").Value = 2323.0588
$ws.Range("I46").Value = 1689.8
$ws.Range("J46").Value = 2586.9167
$ws.Range("K46").Value = 1689.8
$ws.Range("L46").Value = 2586.9167
$ws.Range("M46").Value = -1501.8
$ws.Range("N46").Value = -2962.9167
# LTW row 59
$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -31308
# LTW row 82
$ws.Range("H82").Value = 3334.95
$ws.Range("I82").Value = 1483.3334
$ws.Range("K82").Value = 1483.3334
$ws.Range("M82").Value = -1122.3334
# LTW row 85
$ws.Range("H85").Value = 3334.95
$ws.Range("I85").Value = 1483.3334
$ws.Range("K85").Value = 1483.3334
$ws.Range("M85").Value = -235.3334
# LTW row 100
$ws.Range("H100").Value = 2395.2
$ws.Range("J100").Value = 2744
$ws.Range("L100").Value = 2744
$ws.Range("N100").Value = -3826

$ws = $wb.Worksheets.Item("WVR")
# WVR row 2
$ws.Range("H2").Value = 56170
$ws.Range("I2").Value = 56950
$ws.Range("J2").Value = 55975
$ws.Range("K2").Value = 56950
$ws.Range("L2").Value = 55975
$ws.Range("M2").Value = -56838
$ws.Range("N2").Value = -56199
